$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D values to remain text (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.208.74"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").Value = "1.851.31"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("E4").Value = "  -0.48%  "

$ws.Range("D5").Value = "313.31"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("D7").Value = "0.4606"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "0.3708"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  -1.00%  "

$ws.Range("D10").Value = "0.8833"
$ws.Range("E10").Value = "  +0.85%  "

$ws.Range("D11").Value = "20.02"
$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07817"
$ws.Range("E12").Value = "  -1.50%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.902.46"
$ws.Range("E13").Value = "  +2.95%  "

$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("D15").Value = "6.497"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("D16").Value = "91.31"
$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "0.000008933"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("D19").Value = "1.002"

$ws.Range("D20").Value = "14.69"
$ws.Range("E20").Value = "  -0.74%  "

$ws.Range("D21").Value = "27.238.67"
$ws.Range("E21").Value = "  +1.27%  "

$ws.Range("D22").Value = "5.070"
$ws.Range("E22").Value = "  -0.87%  "

$ws.Range("E23").Value = "  -0.82%  "

$ws.Range("D24").Value = "2.140.92"
$ws.Range("E24").Value = "  +6.74%  "

$ws.Range("D25").Value = "1.952"
$ws.Range("E25").Value = "  +6.04%  "

$ws.Range("D26").Value = "151.91"
$ws.Range("E26").Value = "  -0.84%  "

$ws.Range("D27").Value = "18.41"
$ws.Range("E27").Value = "  -0.13%  "

$ws.Range("D28").Value = "2.061"
$ws.Range("E28").Value = "  +0.71%  "

$ws.Range("D29").Value = "115.57"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("D30").Value = "5.043"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").Value = "0.08823"
$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("D32").Value = "3.092"
$ws.Range("E32").Value = "  +4.30%  "

$ws.Range("D33").Value = "0.7627"
$ws.Range("E33").Value = "  +3.96%  "

$ws.Range("D34").Value = "1.168"
$ws.Range("E34").Value = "  +3.32%  "

$ws.Range("D35").Value = "4.498"

$ws.Range("D36").Value = "2.730"
$ws.Range("E36").Value = "  +9.93%  "

$ws.Range("D37").Value = "1.082"
$ws.Range("E37").Value = "  +0.94%  "

$ws.Range("E38").Value = "  -0.62%  "

$ws.Range("D39").Value = "0.05231"
$ws.Range("E39").Value = "  -0.14%  "

$ws.Range("D40").Value = "2.941"

$ws.Range("D41").Value = "7.068"
$ws.Range("E41").Value = "  -0.54%  "

$ws.Range("E42").Value = "  -1.32%  "

$ws.Range("D43").Value = "0.1624"
$ws.Range("E43").Value = "  -0.31%  "

$ws.Range("D44").Value = "8.377"
$ws.Range("E44").Value = "  +2.07%  "

$ws.Range("D45").Value = "0.4780"
$ws.Range("E45").Value = "  -1.39%  "

$ws.Range("D46").Value = "10.28"
$ws.Range("E46").Value = "  +0.99%  "

$ws.Range("E47").Value = "  -0.50%  "

$ws.Range("E48").Value = "  +0.31%  "

$ws.Range("D49").Value = "1.634"
$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("E50").Value = "  +0.21%  "

$ws.Range("D51").Value = "65.76"
$ws.Range("E51").Value = "  +1.30%  "

# Restore default (unstyled) formatting for column D after forcing text
$ws.Range("D2:D51").Style = "Normal"
